$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 3 (place 1, 100m dash): clear Athlete/School/Time (remove the runs
# entirely, leaving empty paragraphs), reset "a" score to 0
$t.Cell(3, 2).Range.Delete()
$t.Cell(3, 3).Range.Delete()
$t.Cell(3, 4).Range.Delete()
$t.Cell(3, 5).Range.Text = "0"

# Row 5 (place 3, 100m dash): clear Athlete/School/Time, reset "b" score to 0
$t.Cell(5, 2).Range.Delete()
$t.Cell(5, 3).Range.Delete()
$t.Cell(5, 4).Range.Delete()
$t.Cell(5, 6).Range.Text = "0"

# Row 6 (Total, 100m dash): reset "a"/"b" totals to 0
$t.Cell(6, 5).Range.Text = "0"
$t.Cell(6, 6).Range.Text = "0"
